# Add two new register rows ("V and I" RMS data and "Energy" data) to the
# registers table on Sheet1, just above the "<! Regular Registers>" marker
# row, and update the sheet view / column width bookkeeping accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "<! Regular Registers>" marker currently lives on row 28, with rows
# 26 and 27 empty. Push it (and anything below it) down by one row so we
# can populate rows 26 and 27 with the new register entries.
$ws.Rows.Item(28).Insert()

# New row 26: "V and I" (RMS voltage/current) register
# New row 27: "Energy" register
# (Labels, then descriptions, then the numeric/access columns - this
# mirrors the order in which the original author populated the cells.)
$ws.Cells.Item(26, 1).Value2 = "V and I"
$ws.Cells.Item(27, 1).Value2 = "Energy"
$ws.Cells.Item(27, 5).Value2 = "2*3*2 Bytes (Whr,VARhr,Vahr) (2 bytes each, two phase each AC and BC)"
$ws.Cells.Item(26, 5).Value2 = "3*2*2 Bytes (VRMS, IRMS ) (3 bytes each, two phase each AC and BC)"

$ws.Cells.Item(26, 2).Value2 = 14
$ws.Cells.Item(26, 3).Value2 = 18
$ws.Cells.Item(26, 4).Value2 = "R"
$ws.Cells.Item(27, 2).Value2 = 15
$ws.Cells.Item(27, 3).Value2 = 18
$ws.Cells.Item(27, 4).Value2 = "R"

# Widen column E to fit the new, longer description text.
$ws.Columns.Item(5).ColumnWidth = 64

# Scroll the view down a bit and move the selection to the marker row,
# mirroring where the author's cursor ended up after the edit.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A28").Select()
